$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$ws.Range("K2").Value = -139.5
$ws.Range("U2").Value = 85.90000000000001
$ws.Range("V2").Value = 0.07141075733643694
$ws.Range("W2").Value = -0.5537911869789599
$ws.Range("X2").Value = 0.06308892237090165
$ws.Range("Y2").Value = -0.6168801093498616
$ws.Range("AA2").Value = -0.7305435932017434
$ws.Range("AB2").Value = 0.06302289274102885
$ws.Range("AC2").Value = -0.7935664859427722
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 2.67172783334651
$ws.Range("AF2").Value = 2.67172783334651
$ws.Range("AG2").Value = -83.2282721666535
$ws.Range("AH2").Value = 0.002216150040402937
$ws.Range("AI2").Value = 0.007397389184845915
$ws.Range("AJ2").Value = -0.07433274423004928
$ws.Range("AK2").Value = -0.3023495105063644
$ws.Range("AM2").Value = -2.23
$ws.Range("AN2").Value = -0
$ws.Range("AP2").Value = 0.626813316513432
$ws.Range("AQ2").Value = 61.43497757847533

$ws.Range("K3").Value = -139.5
$ws.Range("U3").Value = 85.90000000000001
$ws.Range("V3").Value = 0.07141075733643694
$ws.Range("W3").Value = -0.5537911869789599
$ws.Range("X3").Value = 0.06308892237090165
$ws.Range("Y3").Value = -0.6168801093498616
$ws.Range("AA3").Value = -0.7305435932017434
$ws.Range("AB3").Value = 0.06302289274102885
$ws.Range("AC3").Value = -0.7935664859427722
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 2.67172783334651
$ws.Range("AF3").Value = 2.67172783334651
$ws.Range("AG3").Value = -83.2282721666535
$ws.Range("AH3").Value = 0.002216150040402937
$ws.Range("AI3").Value = 0.007397389184845915
$ws.Range("AJ3").Value = -0.07433274423004928
$ws.Range("AK3").Value = -0.3023495105063644
$ws.Range("AM3").Value = -2.23
$ws.Range("AN3").Value = -0
$ws.Range("AP3").Value = 0.626813316513432
$ws.Range("AQ3").Value = 61.43497757847533

